$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Nathan May): score increases, total recalculated (bonus folded in)
$ws.Range("C2").Value = 18
$ws.Range("F2").Value = 20

# Rows 3 & 4 swap order: Yannick den Daggelder now ranks above Louis Tweddle
$ws.Range("B3").Value = "Yannick den Daggelder"
$ws.Range("C3").Value = 17
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 18

$ws.Range("B4").Value = "Louis Tweddle"
$ws.Range("C4").Value = 17
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 18

# Row 7 (Milan Schoenmakers): score/total bumped by folded-in bonus
$ws.Range("C7").Value = 11
$ws.Range("F7").Value = 12

# Row 11 (Alessandro Delia): score/total bumped by folded-in bonus
$ws.Range("C11").Value = 8
$ws.Range("F11").Value = 8

# Rows 12 & 13 swap order: Diego Meerveld now ranks above Magnus Gladh
$ws.Range("B12").Value = "Diego Meerveld"
$ws.Range("C12").Value = 4
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 4

$ws.Range("B13").Value = "Magnus Gladh"
$ws.Range("C13").Value = 4
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 4

# Row 15 (Robin Willis): score/total bumped by folded-in bonus
$ws.Range("C15").Value = 2
$ws.Range("F15").Value = 2

# Remove the "Extra punten" / "Totaal inclusief bonus" columns entirely;
# their values have already been folded into Score/Totaal above.
$ws.Range("G1:H18").EntireColumn.Delete()
